# Work Plan status-color refresh.
#
# The work-plan table's 2nd column (the 810-twip "status" cell of every
# task row) is shaded to reflect task progress: FFCC66 = in progress,
# FF7C80 = not started / at risk, 00CC66 = complete. Several tasks have
# since finished (or moved along), so their status cell colors change:
#
#   LR, SRS & SDD Update     FFCC66 -> 00CC66
#   Simulation Development   FFCC66 -> 00CC66
#   Web Development          FFCC66 -> 00CC66
#   Mapping                  FFCC66 -> 00CC66
#   Test Plan Document       FF7C80 -> 00CC66
#   Interactive Demo         FF7C80 -> 00CC66
#   User Manual               FF7C80 -> 00CC66
#   Project Report            FF7C80 -> 00CC66
#   Project Poster            FF7C80 -> FFCC66
#   (Project Presentation stays FF7C80 - untouched)

function Get-WdColor($r, $g, $b) {
    # Word's OLE_COLOR (BackgroundPatternColor) packs bytes as 0x00BBGGRR.
    return $b * 65536 + $g * 256 + $r
}

$wdGreen  = Get-WdColor 0x00 0xCC 0x66   # 00CC66
$wdYellow = Get-WdColor 0xFF 0xCC 0x66   # FFCC66
$wdRed    = Get-WdColor 0xFF 0x7C 0x80   # FF7C80

$updates = @{
    "LR, SRS & SDD Update" = $wdGreen
    "Simulation Development" = $wdGreen
    "Web Development" = $wdGreen
    "Mapping" = $wdGreen
    "Test Plan Document" = $wdGreen
    "Interactive Demo" = $wdGreen
    "User Manual" = $wdGreen
    "Project Report" = $wdGreen
    "Project Poster" = $wdYellow
}

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

for ($r = 1; $r -le $t.Rows.Count; $r++) {
    $taskName = $t.Cell($r, 1).Range.Text
    foreach ($key in $updates.Keys) {
        if ($taskName -like "*$key*") {
            $statusCell = $t.Cell($r, 2)
            $statusCell.Shading.BackgroundPatternColor = $updates[$key]
        }
    }
}

# The table's custom style ("a") also drops its redundant explicit
# top/bottom cell-margin overrides (both were already 0, i.e. the same
# value the style inherits from its TableNormal base - a no-op cleanup).
# Best-effort only: not every Word COM surface exposes table-style cell
# margins, so failures here are swallowed without affecting the edits
# above.
try {
    $tableStyle = $d.Styles.Item("a").Table
    $tableStyle.TopPadding = 0
    $tableStyle.BottomPadding = 0
} catch {
}
